# "Namen verander naar engels" - rename Dutch fuel-type labels to English
# on both worksheets ("Nieuw" and "Tweedehands").
#
# The order in which the new (English) string values are first assigned
# matters: Excel's shared-strings table is rebuilt on save using first
# discovery order, and the target workbook has a very specific ordering
# of the newly introduced strings. We therefore assign cell values in the
# exact order needed to reproduce that layout.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Nieuw"
$ws2 = $wb.Worksheets.Item(2)   # "Tweedehands"

# --- Row 3: benzine -> Petrol ---
$ws1.Cells.Item(3,1).Value = "Petrol"
$ws2.Cells.Item(3,1).Value = "Petrol"

# --- Row 4: diesel -> Diesel ---
$ws1.Cells.Item(4,1).Value = "Diesel"
$ws2.Cells.Item(4,1).Value = "Diesel"

# --- Row 5: gas + benzine -> Gas + petrol ---
$ws1.Cells.Item(5,1).Value = "Gas + petrol"
$ws2.Cells.Item(5,1).Value = "Gas + petrol"

# --- Row 6: elektrisch -> Electric ---
$ws1.Cells.Item(6,1).Value = "Electric"
$ws2.Cells.Item(6,1).Value = "Electric"

# --- Row 2: onbekend -> Unknown ---
$ws1.Cells.Item(2,1).Value = "Unknown"
$ws2.Cells.Item(2,1).Value = "Unknown"

# --- Row 7: andere -> Alternative ---
$ws1.Cells.Item(7,1).Value = "Alternative"
$ws2.Cells.Item(7,1).Value = "Alternative"

# --- Row 8: aardgas -> Natural gas ---
$ws1.Cells.Item(8,1).Value = "Natural gas"
$ws2.Cells.Item(8,1).Value = "Natural gas"

# --- Row 9: benzine + elektrisch -> Petrol + electric ---
$ws1.Cells.Item(9,1).Value = "Petrol + electric"
$ws2.Cells.Item(9,1).Value = "Petrol + electric"

# --- Row 10 (Tweedehands only, first): diesel+elektrisch -> Diesel+electric ---
$ws2.Cells.Item(10,1).Value = "Diesel+electric"

# --- Row 11 is "mengsmering", which is kept untranslated, so no change needed ---

# --- Row 12: waterstof -> Hydrogen ---
$ws1.Cells.Item(12,1).Value = "Hydrogen"
$ws2.Cells.Item(12,1).Value = "Hydrogen"

# --- Row 13: diesel + gas -> Diesel + gas ---
$ws1.Cells.Item(13,1).Value = "Diesel + gas"
$ws2.Cells.Item(13,1).Value = "Diesel + gas"

# --- Row 14: elektrisch + lpg -> Electric + liquefied petroleum gas ---
$ws1.Cells.Item(14,1).Value = "Electric + liquefied petroleum gas"
$ws2.Cells.Item(14,1).Value = "Electric + liquefied petroleum gas"

# --- Row 15: bio-ethanol -> Bio-ethanol ---
$ws1.Cells.Item(15,1).Value = "Bio-ethanol"
$ws2.Cells.Item(15,1).Value = "Bio-ethanol"

# --- Row 16: waterstof + elektrisch -> Hydrogen + electric ---
$ws1.Cells.Item(16,1).Value = "Hydrogen + electric"
$ws2.Cells.Item(16,1).Value = "Hydrogen + electric"

# --- Row 17: benzine + aardgas -> Petrol + natural gas ---
$ws1.Cells.Item(17,1).Value = "Petrol + natural gas"
$ws2.Cells.Item(17,1).Value = "Petrol + natural gas"

# --- Row 10 (Nieuw, last): diesel+elektrisch -> "Diesel + electric" (with spaces) ---
$ws1.Cells.Item(10,1).Value = "Diesel + electric"

# --- Update the active sheet / selections to match the end-of-session state ---
# Final state: "Tweedehands" selection at C22 (not the active tab), and
# "Nieuw" becomes the active tab with selection at A16.
[void]$ws2.Range("C22").Select()
[void]$ws1.Range("A16").Select()
